$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'58.497.59"
$ws.Range("E2").Value = "'  -4.01%  "
$ws.Range("D3").Value = "'3.192.49"
$ws.Range("E3").Value = "'  -4.86%  "
$ws.Range("E4").Value = "'  +0.10%  "
$ws.Range("D5").Value = "'533.61"
$ws.Range("E5").Value = "'  -5.79%  "
$ws.Range("D6").Value = "'134.78"
$ws.Range("E6").Value = "'  -8.38%  "
$ws.Range("E7").Value = "'  -0.07%  "
$ws.Range("D8").Value = "'3.190.69"
$ws.Range("E8").Value = "'  -4.88%  "
$ws.Range("E9").Value = "'  -4.93%  "
$ws.Range("D10").Value = "'7.42"
$ws.Range("E10").Value = "'  -6.60%  "
$ws.Range("E11").Value = "'  -7.26%  "
$ws.Range("D12").Value = "'0.394"
$ws.Range("E12").Value = "'  -5.21%  "
$ws.Range("D13").Value = "'3.739.62"
$ws.Range("E13").Value = "'  -5.09%  "
$ws.Range("E14").Value = "'  -0.56%  "
$ws.Range("E15").Value = "'  -8.04%  "
$ws.Range("D16").Value = "'3.194.81"
$ws.Range("E16").Value = "'  -5.19%  "
$ws.Range("D17").Value = "'58.621.48"
$ws.Range("E18").Value = "'  -7.73%  "
$ws.Range("E19").Value = "'  -6.59%  "
$ws.Range("D20").Value = "'13.19"
$ws.Range("E20").Value = "'  -8.55%  "
$ws.Range("D21").Value = "'8.11"
$ws.Range("E21").Value = "'  -9.18%  "
$ws.Range("D22").Value = "'358.46"
$ws.Range("E22").Value = "'  -4.79%  "
$ws.Range("E23").Value = "'  +0.01%  "
$ws.Range("D24").Value = "'69.67"
$ws.Range("E24").Value = "'  -7.03%  "
$ws.Range("D25").Value = "'0.516"
$ws.Range("E25").Value = "'  -7.95%  "
$ws.Range("D26").Value = "'3.326.85"
$ws.Range("E26").Value = "'  -5.10%  "
$ws.Range("D27").Value = "'0.169"
$ws.Range("E27").Value = "'  -3.46%  "
$ws.Range("D28").Value = "'0.0₃0948"
$ws.Range("E28").Value = "'  -12.14%  "
$ws.Range("E29").Value = "'  -1.13%  "
$ws.Range("D30").Value = "'7.04"
$ws.Range("E30").Value = "'  -4.82%  "
$ws.Range("E31").Value = "'  +0.01%  "
$ws.Range("E32").Value = "'  -8.10%  "
$ws.Range("D33").Value = "'7.01"
$ws.Range("E33").Value = "'  -8.87%  "
$ws.Range("D34").Value = "'21.64"
$ws.Range("E34").Value = "'  -5.38%  "
$ws.Range("E35").Value = "'  -8.58%  "
$ws.Range("B36").Value = "'NEARProtocol"
$ws.Range("C36").Value = "'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D36").Value = "'4.94"
$ws.Range("E36").Value = "'  -7.42%  "
$ws.Range("B37").Value = "'Monero"
$ws.Range("C37").Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D37").Value = "'160.87"
$ws.Range("E37").Value = "'  -5.00%  "
$ws.Range("D38").Value = "'6.33"
$ws.Range("E38").Value = "'  -6.96%  "
$ws.Range("E39").Value = "'  -8.60%  "
$ws.Range("D40").Value = "'25.75"
$ws.Range("E40").Value = "'  -10.80%  "
$ws.Range("D41").Value = "'0.0703"
$ws.Range("E41").Value = "'  -6.41%  "
$ws.Range("D42").Value = "'3.220.94"
$ws.Range("E42").Value = "'  -5.05%  "
$ws.Range("D43").Value = "'40.75"
$ws.Range("E43").Value = "'  -3.68%  "
$ws.Range("D44").Value = "'0.708"
$ws.Range("E44").Value = "'  -6.69%  "
$ws.Range("E45").Value = "'  -3.65%  "
$ws.Range("D46").Value = "'4.01"
$ws.Range("E46").Value = "'  -6.73%  "
$ws.Range("B47").Value = "'FirstDigitalUSD"
$ws.Range("C47").Value = "'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D47").Value = "'1.00"
$ws.Range("E47").Value = "'  +0.04%  "
$ws.Range("B48").Value = "'Stacks"
$ws.Range("C48").Value = "'https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D48").Value = "'1.48"
$ws.Range("E48").Value = "'  -7.75%  "
$ws.Range("D49").Value = "'2.282.52"
$ws.Range("E49").Value = "'  -8.44%  "
$ws.Range("D50").Value = "'6.25"
$ws.Range("E50").Value = "'  -6.43%  "
$ws.Range("D51").Value = "'20.42"
$ws.Range("E51").Value = "'  -9.84%  "
